$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-11 (A: id, B: nome, C: departamento, D: motivo, E: horas, F: data serial, G: salario)
$data = @(
    @(38314, "Isadora Silveira", "P&D", "Problemas pessoais", 6, 45086, 3252.39),
    @(83976, "Juan Novaes", "Marketing", "Doença", 5, 45105, 9931.33),
    @(43694, "Milena Cunha", "Jurídico", "Doença", 8, 45094, 6917.57),
    @(69753, "Lorenzo Costa", "Jurídico", "Outros", 7, 45101, 9467.85),
    @(50093, "Dra. Bárbara Pinto", "Engenharia", "Doença", 6, 45094, 2550.47),
    @(11928, "Nina Almeida", "Operações", "Consulta médica", 1, 45084, 4861.24),
    @(56997, "Catarina Almeida", "Financeiro", "Problemas pessoais", 5, 45084, 11975.93),
    @(31210, "Sabrina Costa", "Jurídico", "Problemas pessoais", 5, 45096, 6892.67),
    @(58736, "Dr. Theo Gonçalves", "Vendas", "Consulta médica", 2, 45093, 11105.83),
    @(1593, "Camila Carvalho", "Atendimento ao Cliente", "Outros", 2, 45089, 11967.64)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $ws.Cells.Item($row, 7).Value = $r[6]
    $row++
}
